# Update the "supervision" CV table with the two new education entries
# (replacing the two previous ones), matching the commit:
# "Versión completa y en línea por GitHub Pages ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: MSc in Psychology / 2019 - 2020 / Yenny Johanna Baron Londoño / Universidad El Bosque
$ws.Range("A3").Value = "MSc in Psychology"
$ws.Range("B3").Value = "2019 - 2020"
$ws.Range("C3").Value = "Yenny Johanna Baron Londoño"
$ws.Range("D3").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"

# Row 2: MSc in Neuropsychology / 2022-2023 / Sara Silva Gómez / Universidad Internacional de Valencia
$ws.Range("A2").Value = "MSc in Neuropsychology"
$ws.Range("B2").Value = "2022-2023"
$ws.Range("C2").Value = "Sara Silva Gómez"
$ws.Range("D2").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

# Row heights adjusted (content length / wrapping changed)
$ws.Rows.Item(2).RowHeight = 57.6
$ws.Rows.Item(3).RowHeight = 74.4

# Selection moved to E2:E3
[void]$ws.Range("E2:E3").Select()
